$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45086
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("R2").Value = 'Región del Maule'
$ws.Range("S2").Value = 1000
$ws.Range("D3").Value = 44344
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13500
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 750
$ws.Range("D4").Value = 45112
$ws.Range("M4").Value = 30
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = '$/caja 12 kilos granel'
$ws.Range("R4").Value = 'Región del Maule'
$ws.Range("S4").Value = 13000
$ws.Range("T4").Value = 1
$ws.Range("D6").Value = 45092
$ws.Range("M6").Value = 35
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18571
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región del Maule'
$ws.Range("S6").Value = 1032
$ws.Range("T6").Value = 18
$ws.Range("D7").Value = 44334
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12500
$ws.Range("Q7").Value = '$/caja 12 kilos empedrada'
$ws.Range("S7").Value = 1042
$ws.Range("T7").Value = 12
$ws.Range("D8").Value = 45113
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 13000
$ws.Range("Q8").Value = '$/caja 12 kilos granel'
$ws.Range("S8").Value = 13000
$ws.Range("T8").Value = 1
$ws.Range("D9").Value = 44698
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 16500
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 917
$ws.Range("D10").Value = 45107
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 1000
$ws.Range("T10").Value = 18
$ws.Range("D11").Value = 45085
$ws.Range("M11").Value = 60
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18500
$ws.Range("S11").Value = 1028
$ws.Range("D12").Value = 45096
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Región del Maule'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 18
$ws.Range("D13").Value = 44316
$ws.Range("N13").Value = 17500
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 17750
$ws.Range("Q13").Value = '$/caja 16 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1109
$ws.Range("T13").Value = 16
$ws.Range("D14").Value = 44316
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 40
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("Q14").Value = '$/caja 16 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 1000
$ws.Range("T14").Value = 16
$ws.Range("D15").Value = 45099
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 18000
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 18
$ws.Range("D16").Value = 45111
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 13000
$ws.Range("O16").Value = 13000
$ws.Range("P16").Value = 13000
$ws.Range("Q16").Value = '$/caja 12 kilos granel'
$ws.Range("S16").Value = 13000
$ws.Range("T16").Value = 1